# HandRankingStat.xlsx update
# - Rebuild hand-ranking table (new hand types / balance numbers, 2 new rows)
# - Re-point active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 header: IncreaseChip / IncreaseDrainage labels (unchanged text,
# but shared-string slot moved around in the source file; values stay the same) ----
$ws.Range("F1").Value = "IncreaseChip"
$ws.Range("G1").Value = "IncreaseDrainage"

# ---- Hand ranking data rows (2-11 replaced, 12-13 newly added) ----
$data = @(
    @(2,  "FlUSH FIVE CARD", 1, 180, 14, 0, 70, 10),
    @(3,  "FIVE CARD",       1, 170, 13, 0, 60,  9),
    @(4,  "ROYAL FLUSH",     1, 100, 12, 0, 50,  8),
    @(5,  "FOUR CARD",       1, 140, 10, 0, 40,  6),
    @(6,  "FULL HOUSE",      1, 160,  4, 0, 35,  4),
    @(7,  "STRAIGHT FLUSH",  1, 120,  4, 0, 35,  4),
    @(8,  "FLUSH",           1, 180,  4, 0, 35,  4),
    @(9,  "STRAIGHT",        1, 200,  4, 0, 30,  4),
    @(10, "TRIPLE",          1, 220,  4, 0, 30,  3),
    @(11, "TWO PAIR",        1, 240,  3, 0, 20,  3),
    @(12, "ONE PAIR",        1,  10,  2, 0, 20,  2),
    @(13, "HIGH CARD",       1,  10,  1, 0, 10,  1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]

    # numeric columns are centered, same as the rest of the table
    $ws.Range("B$($r):G$($r)").HorizontalAlignment = -4108
    $ws.Range("B$($r):G$($r)").VerticalAlignment = -4108
}

# Rows 2 & 3 also carry an explicit centered style on column A (rest of column A
# relies on the column-level style already applied to column A)
$ws.Range("A2:A3").HorizontalAlignment = -4108
$ws.Range("A2:A3").VerticalAlignment = -4108

# ---- Selection moved to A8 ----
$ws.Range("A8").Select() | Out-Null
